$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# cryptos list refresh (price / 1h volume updates, plus a few ranking swaps)

$ws.Range("D2").Value = "84.816.20"
$ws.Range("E2").Value = "  +5.51%  "
$ws.Range("D3").Value = "3.315.36"
$ws.Range("E3").Value = "  +2.60%  "
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'219.96"
$ws.Range("E5").Value = "  +2.76%  "
$ws.Range("D6").Value = "'636.95"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("D7").Value = "'0.323"
$ws.Range("E7").Value = "  +17.50%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -2.78%  "
$ws.Range("D10").Value = "3.328.47"
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("D11").Value = "'0.598"
$ws.Range("E11").Value = "  -4.43%  "
$ws.Range("D12").Value = "'0.0000278"
$ws.Range("E12").Value = "  +0.56%  "
$ws.Range("E13").Value = "  +0.01%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "3.917.69"
$ws.Range("E14").Value = "  +2.50%  "
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").Value = "'34.25"
$ws.Range("E15").Value = "  +3.99%  "
$ws.Range("E16").Value = "  -1.09%  "
$ws.Range("D17").Value = "84.589.94"
$ws.Range("E17").Value = "  +5.56%  "
$ws.Range("D18").Value = "3.309.29"
$ws.Range("E18").Value = "  +1.78%  "
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").Value = "'14.68"
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D20").Value = "'3.21"
$ws.Range("E20").Value = "  +4.44%  "
$ws.Range("D21").Value = "'9.22"
$ws.Range("E21").Value = "  -2.32%  "
$ws.Range("D22").Value = "'439.10"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").Value = "'5.25"
$ws.Range("E23").Value = "  -3.37%  "
$ws.Range("D24").Value = "'7.39"
$ws.Range("E24").Value = "  +2.72%  "
$ws.Range("D25").Value = "'5.46"
$ws.Range("E25").Value = "  +11.46%  "
$ws.Range("D26").Value = "'12.19"
$ws.Range("E26").Value = "  +10.23%  "
$ws.Range("D27").Value = "3.459.15"
$ws.Range("E27").Value = "  +2.73%  "
$ws.Range("D28").Value = "'78.03"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("E29").Value = "  +2.04%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "'604.73"
$ws.Range("E31").Value = "  +5.52%  "
$ws.Range("D32").Value = "'0.164"
$ws.Range("E32").Value = "  +32.48%  "
$ws.Range("D33").Value = "'9.29"
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("D34").Value = "'1.58"
$ws.Range("E34").Value = "  +2.57%  "
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("B36").Value = "PancakeSwap"
$ws.Range("C36").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D36").Value = "'2.04"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").Value = "'0.153"
$ws.Range("E37").Value = "  -3.82%  "
$ws.Range("D38").Value = "'23.32"
$ws.Range("E38").Value = "  -1.88%  "
$ws.Range("D39").Value = "'6.43"
$ws.Range("E39").Value = "  +9.02%  "
$ws.Range("E40").Value = "  +0.16%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  +11.53%  "
$ws.Range("E43").Value = "  +10.38%  "
$ws.Range("E44").Value = "  +3.10%  "
$ws.Range("D45").Value = "'159.43"
$ws.Range("E45").Value = "  -2.87%  "
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'190.82"
$ws.Range("E47").Value = "  -1.80%  "
$ws.Range("D48").Value = "'45.26"
$ws.Range("E48").Value = "  +3.87%  "
$ws.Range("E49").Value = "  -0.44%  "
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").Value = "'26.68"
$ws.Range("E51").Value = "  +2.21%  "
